$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.746.16"
$ws.Range("E2").Value = "  +0.58%  "

$ws.Range("D3").Value = "3.097.50"
$ws.Range("E3").Value = "  +3.82%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "388.92"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.83%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "103.77"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.06%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.544"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("E9").Value = "  -0.43%  "

$ws.Range("E10").Value = "  +1.38%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.138"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.22%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0862"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.49%  "

$ws.Range("D13").Value = "3.588.26"
$ws.Range("E13").Value = "  +3.78%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.71"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.64%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.83"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.70%  "

$ws.Range("D16").Value = "3.097.78"
$ws.Range("E16").Value = "  +3.77%  "

$ws.Range("E17").Value = "  -1.28%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.68"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.45%  "

$ws.Range("D19").Value = "51.834.97"
$ws.Range("E19").Value = "  +0.69%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.19"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.83%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.49"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.80%  "

$ws.Range("E22").Value = "  +0.84%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.06"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.24%  "

$ws.Range("E24").Value = "  +0.64%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.15"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.35%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.16"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.06%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "27.23"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.55%  "

$ws.Range("E28").Value = "  +1.60%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.23"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.29%  "

$ws.Range("E31").Value = "  -0.35%  "

$ws.Range("E32").Value = "  +0.09%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "35.63"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.57%  "

$ws.Range("E34").Value = "  +0.75%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "50.39"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.90%  "

$ws.Range("E36").Value = "  +1.57%  "

$ws.Range("E37").Value = "  -0.20%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.39"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.43%  "

$ws.Range("E39").Value = "  +9.02%  "

$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.89"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.26%  "

$ws.Range("B41").Value = "Celestia"
$ws.Range("C41").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "17.01"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.42%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.58"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.83%  "

$ws.Range("B43").Value = "Stellar"
$ws.Range("C43").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.116"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.38%  "

$ws.Range("B44").Value = "Monero"
$ws.Range("C44").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "127.32"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.92%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.70"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.36%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.17"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.76%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.45"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.74%  "

$ws.Range("E48").Value = "  +2.15%  "

$ws.Range("D49").Value = "2.051.40"
$ws.Range("E49").Value = "  +1.28%  "

$ws.Range("D50").Value = "3.407.26"
$ws.Range("E50").Value = "  +3.77%  "

$ws.Range("E51").Value = "  +7.20%  "
